# Workbook edit script
# - Rename sheet "R=16K" -> "R"
# - Rename sheet "Sheet1" -> "b) Error"
# - Sheet "R": delete the stray row 4 (anomalous frequency 65010 row), clear G1
#   stray cell, and move selection to D3
# - Sheet "b) Error": widen column E, add AVERAGE formula in E10, move
#   selection to D10

$wb = $excel.ActiveWorkbook

$wsR = $wb.Worksheets.Item(1)
$wsError = $wb.Worksheets.Item(4)

# Rename sheets
$wsR.Name = "R"
$wsError.Name = "b) Error"

# --- Sheet "R" edits ---
# Remove the orphan/duplicate "R2[Ohm]" label in G1
$wsR.Range("G1").ClearContents()

# Delete the anomalous row (C4 = 65010, with blank D4/E4) - this shifts
# everything below up by one row and drops the final blank row (47)
$wsR.Rows.Item(4).Delete()

# Update selection to D3
$wsR.Activate()
$wsR.Range("D3").Select()

# --- Sheet "b) Error" edits ---
# Widen column E (target raw xlsx width = 18 => ColumnWidth = 18 - 5/6)
$wsError.Columns.Item(5).ColumnWidth = 18 - 5/6

# Add the average of the error column
$wsError.Range("E10").NumberFormat = "0.00"
$wsError.Range("E10").Formula = "=AVERAGE(E2:E9)"

# Update selection to D10
$wsError.Activate()
$wsError.Range("D10").Select()
